# Scrum.xlsx update: "Finished mvc doc and filled table in scrum file."
#
# - On the "Burndown Chart" sheet, log 9 hours against the "22.04" build
#   task on the 28.04 day column (G8), which is the missing data point
#   feeding the burndown formulas in row 16 and the chart series cache.
# - Make "Burndown Chart" the active/selected sheet (it was "01.05"
#   before), with the view scrolled down to show the newly entered row
#   and the matching cell selected.

$wb = $excel.ActiveWorkbook

$burndown = $wb.Worksheets.Item("Burndown Chart")
$lastDay  = $wb.Worksheets.Item("01.05")

# Fill in the previously-empty data point: 9 hours done on 28.04.
$burndown.Range("G8").Value = 9

# Move focus to the Burndown Chart tab and park the selection/scroll
# position where the new entry lives.
$burndown.Activate()
$burndown.Range("G9").Select()
$burndown.Application.ActiveWindow.ScrollRow = 15

# "01.05" was the previously active tab; it no longer is once Burndown
# Chart becomes active above, matching the dropped tabSelected flag.
